{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"91\u00f78=11, 3\" },\n  { row: 0, col: 1, text: \"95\u00f72=47, 1\" },\n  { row: 0, col: 2, text: \"15\u00f76=2, 3\" },\n  { row: 0, col: 3, text: \"50\u00f79=5, 5\" },\n  { row: 0, col: 4, text: \"88\u00f76=14, 4\" },\n  { row: 4, col: 0, text: \"16\u00f73=5, 1\" },\n  { row: 4, col: 1, text: \"40\u00f76=6, 4\" },\n  { row: 4, col: 2, text: \"42\u00f74=10, 2\" },\n  { row: 4, col: 3, text: \"35\u00f76=5, 5\" },\n  { row: 4, col: 4, text: \"67\u00f72=33, 1\" },\n  { row: 8, col: 0, text: \"56\u00f78=7, 0\" },\n  { row: 8, col: 1, text: \"26\u00f73=8, 2\" },\n  { row: 8, col: 2, text: \"23\u00f75=4, 3\" },\n  { row: 8, col: 3, text: \"26\u00f77=3, 5\" },\n  { row: 8, col: 4, text: \"12\u00f75=2, 2\" },\n  { row: 12, col: 0, text: \"42\u00f76=7, 0\" },\n  { row: 12, col: 1, text: \"28\u00f78=3, 4\" },\n  { row: 12, col: 2, text: \"49\u00f76=8, 1\" },\n  { row: 12, col: 3, text: \"34\u00f72=17, 0\" },\n  { row: 12, col: 4, text: \"48\u00f75=9, 3\" },\n  { row: 16, col: 0, text: \"24\u00f75=4, 4\" },\n  { row: 16, col: 1, text: \"86\u00f79=9, 5\" },\n  { row: 16, col: 2, text: \"21\u00f77=3, 0\" },\n  { row: 16, col: 3, text: \"94\u00f73=31, 1\" },\n  { row: 16, col: 4, text: \"56\u00f79=6, 2\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.body.getRange(\"Whole\").insertText(u.text, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text = '91\u00f78=11, 3'\n$t.Cell(1, 2).Range.Text = '95\u00f72=47, 1'\n$t.Cell(1, 3).Range.Text = '15\u00f76=2, 3'\n$t.Cell(1, 4).Range.Text = '50\u00f79=5, 5'\n$t.Cell(1, 5).Range.Text = '88\u00f76=14, 4'\n$t.Cell(5, 1).Range.Text = '16\u00f73=5, 1'\n$t.Cell(5, 2).Range.Text = '40\u00f76=6, 4'\n$t.Cell(5, 3).Range.Text = '42\u00f74=10, 2'\n$t.Cell(5, 4).Range.Text = '35\u00f76=5, 5'\n$t.Cell(5, 5).Range.Text = '67\u00f72=33, 1'\n$t.Cell(9, 1).Range.Text = '56\u00f78=7, 0'\n$t.Cell(9, 2).Range.Text = '26\u00f73=8, 2'\n$t.Cell(9, 3).Range.Text = '23\u00f75=4, 3'\n$t.Cell(9, 4).Range.Text = '26\u00f77=3, 5'\n$t.Cell(9, 5).Range.Text = '12\u00f75=2, 2'\n$t.Cell(13, 1).Range.Text = '42\u00f76=7, 0'\n$t.Cell(13, 2).Range.Text = '28\u00f78=3, 4'\n$t.Cell(13, 3).Range.Text = '49\u00f76=8, 1'\n$t.Cell(13, 4).Range.Text = '34\u00f72=17, 0'\n$t.Cell(13, 5).Range.Text = '48\u00f75=9, 3'\n$t.Cell(17, 1).Range.Text = '24\u00f75=4, 4'\n$t.Cell(17, 2).Range.Text = '86\u00f79=9, 5'\n$t.Cell(17, 3).Range.Text = '21\u00f77=3, 0'\n$t.Cell(17, 4).Range.Text = '94\u00f73=31, 1'\n$t.Cell(17, 5).Range.Text = '56\u00f79=6, 2'\n"}
